$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: bet amount bumped, add Hit?/Net Result columns ---
$ws.Range("C2").Value = 5
$ws.Range("G2").Value = "yes"
$ws.Range("H2").Formula = "=IF(G2=""Yes"",E2*C2-C2,-C2)"

# --- Row 3: add Hit?/Net Result columns ---
$ws.Range("G3").Value = "no"
$ws.Range("H3").Formula = "=IF(G3=""Yes"",E3*C3-C3,-C3)"

# --- New rows 4-9: additional bet legs ---
$ws.Range("A4").Value = 45382
$ws.Range("A4").NumberFormat = "d-mmm"
$ws.Range("B4").Value = "Garland 4+ reb"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2.71
$ws.Range("E4").Value = 3.6
$ws.Range("F4").Formula = "=E4/D4-1"

$ws.Range("A5").Value = 45382
$ws.Range("A5").NumberFormat = "d-mmm"
$ws.Range("B5").Value = "Garland 10+ pts"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1.02
$ws.Range("E5").Value = 1.12
$ws.Range("F5").Formula = "=E5/D5-1"

$ws.Range("A6").Value = 45382
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("B6").Value = "mobley 10+ pts"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1.03
$ws.Range("E6").Value = 1.24
$ws.Range("F6").Formula = "=E6/D6-1"

$ws.Range("A7").Value = 45382
$ws.Range("A7").NumberFormat = "d-mmm"
$ws.Range("B7").Value = "strus 10+ pts"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1.25
$ws.Range("E7").Value = 1.83
$ws.Range("F7").Formula = "=E7/D7-1"

$ws.Range("A8").Value = 45382
$ws.Range("A8").NumberFormat = "d-mmm"
$ws.Range("B8").Value = "KCP 10+ pts"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1.6
$ws.Range("E8").Value = 2.35
$ws.Range("F8").Formula = "=E8/D8-1"

$ws.Range("A9").Value = 45382
$ws.Range("A9").NumberFormat = "d-mmm"
$ws.Range("B9").Value = "Gordon 10+ pts"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1.08
$ws.Range("E9").Value = 1.23
$ws.Range("F9").Formula = "=E9/D9-1"

# --- Selection moves to F17 ---
$ws.Range("F17").Select()

# --- Best-effort: shrink the saved window chrome to match the author's
#     resized Excel window (not all COM hosts persist this to the
#     workbookView, so failures here are harmless and ignored). ---
try {
    $excel.ActiveWindow.Width = 14400
    $excel.ActiveWindow.Height = 15600
    $excel.ActiveWindow.Left = 0
    $excel.ActiveWindow.Top = 0
} catch {}
